# Added Update/Delete functions in device42 project
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the header row labels (row 1, columns A:O)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "type"
$ws.Range("D1").Value = "network_device"
$ws.Range("E1").Value = "blade_chassis"
$ws.Range("F1").Value = "virtual_host"
$ws.Range("G1").Value = "in_service"
$ws.Range("H1").Value = "service_level"
$ws.Range("I1").Value = "serial_no"
$ws.Range("J1").Value = "asset_no"
$ws.Range("K1").Value = "uuid"
$ws.Range("L1").Value = "customer"
$ws.Range("M1").Value = "blade_slot_no"
$ws.Range("N1").Value = "blade_slot_no"
$ws.Range("O1").Value = "device_host_chassis"

# Set column B width (stored width of 20 corresponds to a ColumnWidth
# COM value of 19 + 1/7 given the default Calibri 11 font metrics)
$ws.Columns.Item(2).ColumnWidth = 19 + (1/7)

# Update selection to P4
$ws.Range("P4").Select()
